# Generate Report for Handback
# Update the "generate date / handoff / handback" timestamp text cells
# on each worksheet, as plain text strings (they are stored as shared
# strings formatted to look like dates, not numeric date values).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: row for 15dfc46c-...md -> "Latest HO Xliff Generate Date" (col G)
$wsOverview.Range("G3").Value = "2016-08-16 18:40:16"

# zh-cn sheet: row for 15dfc46c-...md
#   Correspond Handoff Datetime (col H)
#   Correspond Handback DateTime (col K)
$wsZhCn.Range("H3").Value = "2016-08-16 18:39:59"
$wsZhCn.Range("K3").Value = "2016-08-16 18:40:33"

# de-de sheet: row for 15dfc46c-...md
#   Correspond Handback DateTime (col K)
$wsDeDe.Range("K3").Value = "2016-08-16 18:40:41"
